$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 21.97750833333333
$ws.Cells.Item(2, 8).Value = 65.932525
$ws.Cells.Item(2, 9).Value = 0.5427578249542736
$ws.Cells.Item(2, 10).Value = 0.5427578249542736
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.155977333333333
$ws.Cells.Item(2, 14).Value = 9.467932
$ws.Cells.Item(2, 15).Value = 0.3579027849973545
$ws.Cells.Item(2, 16).Value = 0.3579027849973545
$ws.Cells.Item(2, 17).Value = 69.36051814314443
$ws.Cells.Item(2, 18).Value = 624.2446632883
$ws.Cells.Item(2, 19).Value = 0.1942545371302411
$ws.Cells.Item(2, 20).Value = 0.1942545371302411

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 21.97750833333333
$ws.Cells.Item(3, 8).Value = 65.932525
$ws.Cells.Item(3, 9).Value = 0.5427578249542736
$ws.Cells.Item(3, 10).Value = 0.5427578249542736
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.165953666666667
$ws.Cells.Item(3, 14).Value = 9.497861
$ws.Cells.Item(3, 15).Value = 0.359034148472735
$ws.Cells.Item(3, 16).Value = 0.359034148472735
$ws.Cells.Item(3, 17).Value = 69.5797730921139
$ws.Cells.Item(3, 18).Value = 626.217957829025
$ws.Cells.Item(3, 19).Value = 0.1948685935093714
$ws.Cells.Item(3, 20).Value = 0.1948685935093714

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 21.97750833333333
$ws.Cells.Item(4, 8).Value = 65.932525
$ws.Cells.Item(4, 9).Value = 0.5427578249542736
$ws.Cells.Item(4, 10).Value = 0.5427578249542736
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 2.496042666666666
$ws.Cells.Item(4, 14).Value = 7.488128
$ws.Cells.Item(4, 15).Value = 0.2830630665299106
$ws.Cells.Item(4, 16).Value = 0.2830630665299106
$ws.Cells.Item(4, 17).Value = 54.85679850702221
$ws.Cells.Item(4, 18).Value = 493.7111865632
$ws.Cells.Item(4, 19).Value = 0.1536346943146611
$ws.Cells.Item(4, 20).Value = 0.1536346943146611

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.688376
$ws.Cells.Item(5, 8).Value = 38.065128
$ws.Cells.Item(5, 9).Value = 0.3133528721960219
$ws.Cells.Item(5, 10).Value = 0.3133528721960219
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.155977333333333
$ws.Cells.Item(5, 14).Value = 9.467932
$ws.Cells.Item(5, 15).Value = 0.3579027849973545
$ws.Cells.Item(5, 16).Value = 0.3579027849973545
$ws.Cells.Item(5, 17).Value = 40.04422705281066
$ws.Cells.Item(5, 18).Value = 360.398043475296
$ws.Cells.Item(5, 19).Value = 0.1121498656458763
$ws.Cells.Item(5, 20).Value = 0.1121498656458763

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.688376
$ws.Cells.Item(6, 8).Value = 38.065128
$ws.Cells.Item(6, 9).Value = 0.3133528721960219
$ws.Cells.Item(6, 10).Value = 0.3133528721960219
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.165953666666667
$ws.Cells.Item(6, 14).Value = 9.497861
$ws.Cells.Item(6, 15).Value = 0.359034148472735
$ws.Cells.Item(6, 16).Value = 0.359034148472735
$ws.Cells.Item(6, 17).Value = 40.17081052124534
$ws.Cells.Item(6, 18).Value = 361.537294691208
$ws.Cells.Item(6, 19).Value = 0.1125043816403845
$ws.Cells.Item(6, 20).Value = 0.1125043816403845

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.688376
$ws.Cells.Item(7, 8).Value = 38.065128
$ws.Cells.Item(7, 9).Value = 0.3133528721960219
$ws.Cells.Item(7, 10).Value = 0.3133528721960219
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.496042666666666
$ws.Cells.Item(7, 14).Value = 7.488128
$ws.Cells.Item(7, 15).Value = 0.2830630665299106
$ws.Cells.Item(7, 16).Value = 0.2830630665299106
$ws.Cells.Item(7, 17).Value = 31.67072786670933
$ws.Cells.Item(7, 18).Value = 285.036550800384
$ws.Cells.Item(7, 19).Value = 0.08869862490976115
$ws.Cells.Item(7, 20).Value = 0.08869862490976115

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.008189666666667
$ws.Cells.Item(8, 8).Value = 6.024569
$ws.Cells.Item(8, 9).Value = 0.04959436889042158
$ws.Cells.Item(8, 10).Value = 0.04959436889042158
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 3.155977333333333
$ws.Cells.Item(8, 14).Value = 9.467932
$ws.Cells.Item(8, 15).Value = 0.3579027849973545
$ws.Cells.Item(8, 16).Value = 0.3579027849973545
$ws.Cells.Item(8, 17).Value = 6.337801069034221
$ws.Cells.Item(8, 18).Value = 57.04020962130799
$ws.Cells.Item(8, 19).Value = 0.01774996274606804
$ws.Cells.Item(8, 20).Value = 0.01774996274606804

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.008189666666667
$ws.Cells.Item(9, 8).Value = 6.024569
$ws.Cells.Item(9, 9).Value = 0.04959436889042158
$ws.Cells.Item(9, 10).Value = 0.04959436889042158
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 3.165953666666667
$ws.Cells.Item(9, 14).Value = 9.497861
$ws.Cells.Item(9, 15).Value = 0.359034148472735
$ws.Cells.Item(9, 16).Value = 0.359034148472735
$ws.Cells.Item(9, 17).Value = 6.357835438545444
$ws.Cells.Item(9, 18).Value = 57.220518946909
$ws.Cells.Item(9, 19).Value = 0.01780607200361521
$ws.Cells.Item(9, 20).Value = 0.01780607200361521

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.008189666666667
$ws.Cells.Item(10, 8).Value = 6.024569
$ws.Cells.Item(10, 9).Value = 0.04959436889042158
$ws.Cells.Item(10, 10).Value = 0.04959436889042158
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.496042666666666
$ws.Cells.Item(10, 14).Value = 7.488128
$ws.Cells.Item(10, 15).Value = 0.2830630665299106
$ws.Cells.Item(10, 16).Value = 0.2830630665299106
$ws.Cells.Item(10, 17).Value = 5.01252709075911
$ws.Cells.Item(10, 18).Value = 45.11274381683199
$ws.Cells.Item(10, 19).Value = 0.01403833414073833
$ws.Cells.Item(10, 20).Value = 0.01403833414073833

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3.818218
$ws.Cells.Item(11, 8).Value = 11.454654
$ws.Cells.Item(11, 9).Value = 0.09429493395928291
$ws.Cells.Item(11, 10).Value = 0.09429493395928291
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.155977333333333
$ws.Cells.Item(11, 14).Value = 9.467932
$ws.Cells.Item(11, 15).Value = 0.3579027849973545
$ws.Cells.Item(11, 16).Value = 0.3579027849973545
$ws.Cells.Item(11, 17).Value = 12.05020946172533
$ws.Cells.Item(11, 18).Value = 108.451885155528
$ws.Cells.Item(11, 19).Value = 0.03374841947516897
$ws.Cells.Item(11, 20).Value = 0.03374841947516897

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 3.818218
$ws.Cells.Item(12, 8).Value = 11.454654
$ws.Cells.Item(12, 9).Value = 0.09429493395928291
$ws.Cells.Item(12, 10).Value = 0.09429493395928291
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 3.165953666666667
$ws.Cells.Item(12, 14).Value = 9.497861
$ws.Cells.Item(12, 15).Value = 0.359034148472735
$ws.Cells.Item(12, 16).Value = 0.359034148472735
$ws.Cells.Item(12, 17).Value = 12.08830127723267
$ws.Cells.Item(12, 18).Value = 108.794711495094
$ws.Cells.Item(12, 19).Value = 0.03385510131936392
$ws.Cells.Item(12, 20).Value = 0.03385510131936392

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 3.818218
$ws.Cells.Item(13, 8).Value = 11.454654
$ws.Cells.Item(13, 9).Value = 0.09429493395928291
$ws.Cells.Item(13, 10).Value = 0.09429493395928291
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 2.496042666666666
$ws.Cells.Item(13, 14).Value = 7.488128
$ws.Cells.Item(13, 15).Value = 0.2830630665299106
$ws.Cells.Item(13, 16).Value = 0.2830630665299106
$ws.Cells.Item(13, 17).Value = 9.530435038634666
$ws.Cells.Item(13, 18).Value = 85.77391534771199
$ws.Cells.Item(13, 19).Value = 0.02669141316475003
$ws.Cells.Item(13, 20).Value = 0.02669141316475003
